$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# New diary entry goes into row 113 (row 114 stays empty, same as before).
# Copy the formatting from the row above (row 112) onto the new row first,
# so date/number styles (A=date style, B=centered number style) match
# exactly instead of minting new style records.
$ws.Range("A112:D112").Copy() | Out-Null
$ws.Range("A113:D113").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new row's values.
$ws.Cells.Item(113, 1).Value = 44592   # 2022-01-31
$ws.Cells.Item(113, 2).Value = 3
$ws.Cells.Item(113, 3).Value = "Profile sivun komponentit muutettu käyttämään Sass:ia ja muokattu vähän ilmettä ja koodia"
$ws.Cells.Item(113, 4).Value = "client"

# Extend the hours-total formula to include the new row.
$ws.Range("B115").Formula = "=SUM(B2:B113)"

# Move the active selection to the newly-filled cell, matching the saved UI state.
$ws.Range("D113").Select()

$wb.Save()
